{"js": "// Apply the \"Added many more features\" edit: refresh the SEO title, the\n// What-we-like / What-we-don't-like bullet lists, the bolded title repeat,\n// and the italic meta-description near the end of the document.\n\nconst replacements = [\n  {\n    find: \"Play Fortune Falcon Wild Respins for Free - Slot Game Review\",\n    replace: \"Play Fortune Falcon Wild Respins Slot for Free\",\n  },\n  {\n    find: \"Attractive betting range from \\u20AC0.01 to \\u20AC100\",\n    replace: \"Wide betting range suitable for all players\",\n  },\n  {\n    find: \"Exciting bonus features with free spins and multipliers\",\n    replace: \"Interesting bonus features with potential for big wins\",\n  },\n  {\n    find: \"Detailed visuals and high-quality graphics\",\n    replace: \"High-quality graphics and detailed animations\",\n  },\n  {\n    find: \"Interesting storyline with the falcon as a central character\",\n    replace: \"Enchanting storyline that immerses players in the game\",\n  },\n  {\n    find: \"Lower theoretical RTP of 94.02%\",\n    replace: \"Slightly lower RTP compared to other slots\",\n  },\n  {\n    find: \"Medium volatility might not be suitable for every player\",\n    replace: \"Limited number of pay lines\",\n  },\n  {\n    find: \"Read our review of Fortune Falcon Wild Respins and play for free. Discover exciting bonus features, detailed visuals and attractive betting range.\",\n    replace: \"Read our review of Fortune Falcon Wild Respins and play this slot game for free. Win big with exciting bonus features.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit: refresh the SEO title, the\n# What-we-like / What-we-don't-like bullet lists, the bolded title repeat,\n# and the italic meta-description near the end of the document.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Fortune Falcon Wild Respins for Free - Slot Game Review\"; Replace = \"Play Fortune Falcon Wild Respins Slot for Free\" },\n    @{ Find = \"Attractive betting range from \u20ac0.01 to \u20ac100\"; Replace = \"Wide betting range suitable for all players\" },\n    @{ Find = \"Exciting bonus features with free spins and multipliers\"; Replace = \"Interesting bonus features with potential for big wins\" },\n    @{ Find = \"Detailed visuals and high-quality graphics\"; Replace = \"High-quality graphics and detailed animations\" },\n    @{ Find = \"Interesting storyline with the falcon as a central character\"; Replace = \"Enchanting storyline that immerses players in the game\" },\n    @{ Find = \"Lower theoretical RTP of 94.02%\"; Replace = \"Slightly lower RTP compared to other slots\" },\n    @{ Find = \"Medium volatility might not be suitable for every player\"; Replace = \"Limited number of pay lines\" },\n    @{ Find = \"Read our review of Fortune Falcon Wild Respins and play for free. Discover exciting bonus features, detailed visuals and attractive betting range.\"; Replace = \"Read our review of Fortune Falcon Wild Respins and play this slot game for free. Win big with exciting bonus features.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
